$d = $word.ActiveDocument

# --- Rename the three heading "_Toc..." bookmarks (Word's Bookmark.Name
# setter doesn't write back in this runtime, so re-create each bookmark
# over its own range under the new name instead). ---

function Rename-Bookmark($doc, $oldName, $newName) {
    $bm = $doc.Bookmarks.Item($oldName)
    $r = $bm.Range
    $bm.Delete()
    $doc.Bookmarks.Add($newName, $r) | Out-Null
}

Rename-Bookmark $d "_Toc1636375498552957588038508" "_Toc16363785084833195660448043"
Rename-Bookmark $d "_Toc16363754985868704453244859" "_Toc16363785085207501824900034"
Rename-Bookmark $d "_Toc16363754986284348924878694" "_Toc16363785085562561309543595"

# --- Update the Heading3 text from "jos jedan samo failed" to "etc" ---
$d.Content.Find.Execute("jos jedan samo failed", $true, $false, $false, $false, $false,
                         $true, 1, $false, "etc", 2)
